$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B19").Value = "gar_du_pa_det_studiet_du_hadde_satt_opp_som_ditt_forstevalg"
